# Applies the crypto-price / volume / ranking refresh described in the commit
# "Updated cryptos list on Tue Aug 13 18:17:53 UTC 2024 with GitHub Actions".
#
# All cells on this sheet are plain text (inline/shared strings) -- including the
# "Price" column, which contains values such as "60.791.11" or "0.580" that Excel
# would otherwise silently reinterpret/round as numbers. Set-CellText forces the
# target cell to a text format just long enough to assign the exact literal string,
# then restores the cell style to Normal so no formatting residue is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-CellText "D2" "60.791.11"
$ws.Range("E2").Value = "  +2.69%  "

# Row 3
Set-CellText "D3" "2.685.71"
$ws.Range("E3").Value = "  +1.74%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
Set-CellText "D5" "521.97"
$ws.Range("E5").Value = "  +1.13%  "

# Row 6
Set-CellText "D6" "147.80"
$ws.Range("E6").Value = "  +1.21%  "

# Row 7
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
Set-CellText "D8" "0.580"
$ws.Range("E8").Value = "  +1.68%  "

# Row 9
Set-CellText "D9" "2.704.85"
$ws.Range("E9").Value = "  +1.43%  "

# Row 10
Set-CellText "D10" "6.41"
$ws.Range("E10").Value = "  -0.44%  "

# Row 11
$ws.Range("E11").Value = "  +0.55%  "

# Row 12
Set-CellText "D12" "0.342"
$ws.Range("E12").Value = "  +0.81%  "

# Row 13
$ws.Range("E13").Value = "  +1.42%  "

# Row 14
Set-CellText "D14" "3.160.35"
$ws.Range("E14").Value = "  +1.93%  "

# Row 15
Set-CellText "D15" "60.800.05"
$ws.Range("E15").Value = "  +2.74%  "

# Row 16
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-CellText "D16" "21.47"
$ws.Range("E16").Value = "  +1.10%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-CellText "D17" "2.769.70"
$ws.Range("E17").Value = "  +4.14%  "

# Row 18
$ws.Range("E18").Value = "  +1.10%  "

# Row 19
Set-CellText "D19" "354.40"
$ws.Range("E19").Value = "  +2.61%  "

# Row 20
Set-CellText "D20" "4.58"
$ws.Range("E20").Value = "  -0.06%  "

# Row 21
Set-CellText "D21" "10.56"
$ws.Range("E21").Value = "  +1.08%  "

# Row 22
Set-CellText "D22" "6.36"
$ws.Range("E22").Value = "  +3.50%  "

# Row 23
Set-CellText "D23" "0.999"
$ws.Range("E23").Value = "  -0.04%  "

# Row 24
Set-CellText "D24" "63.08"
$ws.Range("E24").Value = "  +2.63%  "

# Row 25
Set-CellText "D25" "0.424"
$ws.Range("E25").Value = "  -0.10%  "

# Row 26
$ws.Range("E26").Value = "  +4.72%  "

# Row 28
Set-CellText "D28" "0.0₃0825"
$ws.Range("E28").Value = "  +1.21%  "

# Row 29
$ws.Range("E29").Value = "  +1.94%  "

# Row 30
Set-CellText "D30" "6.86"
$ws.Range("E30").Value = "  +5.32%  "

# Row 31
$ws.Range("E31").Value = "  +0.14%  "

# Row 32
Set-CellText "D32" "19.18"
$ws.Range("E32").Value = "  +0.77%  "

# Row 33
$ws.Range("E33").Value = "  +0.88%  "

# Row 34
Set-CellText "D34" "149.58"
$ws.Range("E34").Value = "  -0.46%  "

# Row 35
Set-CellText "D35" "4.27"
$ws.Range("E35").Value = "  +4.41%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-CellText "D36" "1.24"
$ws.Range("E36").Value = "  +6.11%  "

# Row 37
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-CellText "D37" "0.950"
$ws.Range("E37").Value = "  -8.63%  "

# Row 38
Set-CellText "D38" "1.57"
$ws.Range("E38").Value = "  +10.50%  "

# Row 39
Set-CellText "D39" "0.875"
$ws.Range("E39").Value = "  +1.06%  "

# Row 40
$ws.Range("E40").Value = "  +0.21%  "

# Row 41
Set-CellText "D41" "3.74"
$ws.Range("E41").Value = "  +0.36%  "

# Row 42
Set-CellText "D42" "285.21"
$ws.Range("E42").Value = "  -0.17%  "

# Row 43
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-CellText "D43" "0.0994"
$ws.Range("E43").Value = "  +0.89%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-CellText "D44" "20.06"
$ws.Range("E44").Value = "  +2.23%  "

# Row 45
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-CellText "D45" "0.614"
$ws.Range("E45").Value = "  -0.45%  "

# Row 46
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-CellText "D46" "0.995"
$ws.Range("E46").Value = "  +0.21%  "

# Row 47
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-CellText "D47" "2.141.88"
$ws.Range("E47").Value = "  +7.94%  "

# Row 48
Set-CellText "D48" "0.0542"
$ws.Range("E48").Value = "  +0.11%  "

# Row 49
Set-CellText "D49" "4.89"
$ws.Range("E49").Value = "  +5.01%  "

# Row 50
$ws.Range("E50").Value = "  +2.30%  "

# Row 51
Set-CellText "D51" "19.18"
$ws.Range("E51").Value = "  +3.61%  "
